$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell C2 with the value "Red Black Tree", matching the wrap-text style used elsewhere
$ws.Range("C2").Value = "Red Black Tree"
$ws.Range("C2").WrapText = $true

# Update the selection to D20 as reflected in the saved view state
$ws.Range("D20").Select()
